$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-like values that are numeric-looking must be forced to remain as text
# (matching the original inline-string cell type) by temporarily applying a
# text number format, setting the value, then clearing the format so no extra
# style index is left behind on the cell.
function Set-TextValue($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue $ws "D2" "247.33"
Set-TextValue $ws "D3" "21.56"
Set-TextValue $ws "D4" "5.286"
Set-TextValue $ws "D5" "0.05592"
Set-TextValue $ws "D6" "3.395"
Set-TextValue $ws "D7" "6.364"
Set-TextValue $ws "D8" "0.8156"
Set-TextValue $ws "D9" "0.9565"
Set-TextValue $ws "D10" "0.1407"
Set-TextValue $ws "D11" "0.07427"
Set-TextValue $ws "D12" "0.03159"
Set-TextValue $ws "D13" "0.03031"
Set-TextValue $ws "D14" "0.09286"
Set-TextValue $ws "D15" "3.561"
Set-TextValue $ws "D16" "0.001587"
Set-TextValue $ws "D17" "0.04714"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws "D18" "0.0005765"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws "D19" "0.006423"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue $ws "D20" "0.005060"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue $ws "D21" "0.001032"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue $ws "D22" "0.0001499"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws "D23" "3.743"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws "D24" "2.117"
$ws.Range("E24").Value = "23BTSETokenBTSE"
Set-TextValue $ws "D25" "0.3253"
Set-TextValue $ws "D28" "0.0003097"
Set-TextValue $ws "D40" "0.03926"
Set-TextValue $ws "D41" "0.007057"
Set-TextValue $ws "D42" "0.1049"
Set-TextValue $ws "D43" "0.003058"
Set-TextValue $ws "D44" "0.007841"
Set-TextValue $ws "D45" "0.00005806"
Set-TextValue $ws "D46" "0.00000000749"
Set-TextValue $ws "D47" "0.0005495"
Set-TextValue $ws "D48" "0.6794"
Set-TextValue $ws "D49" "0.1557"
$ws.Range("E49").Value = "48BOLOBOLOBestin24h"
Set-TextValue $ws "D50" "0.00002098"
Set-TextValue $ws "D51" "0.01009"
